$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("Rua Monte Sião ","Cidade de Deus ","Rio de Janeiro/RJ ","22770-370")
  ,@("Rua Moreira Campos ","Pechincha ","Rio de Janeiro/RJ ","22770-125")
  ,@("Rua Oscar Cordeiro ","Pechincha ","Rio de Janeiro/RJ ","22770-200")
  ,@("Rua Oswaldo Lussac ","Taquara ","Rio de Janeiro/RJ ","22770-640")
  ,@("Rua Paço do Lumiar ","Pechincha ","Rio de Janeiro/RJ ","22770-080")
  ,@("Rua Paulo Moreira da Silva ","Taquara ","Rio de Janeiro/RJ ","22770-210")
  ,@("Rua Pedro Ivo ","Cidade de Deus ","Rio de Janeiro/RJ ","22770-665")
  ,@("Rua Percy Lau ","Pechincha ","Rio de Janeiro/RJ ","22770-220")
  ,@("Rua Pintor Leandro Joaquim ","Cidade de Deus ","Rio de Janeiro/RJ ","22770-650")
  ,@("Rua Professor Henrique Costa - até 685 - lado ímpar ","Pechincha ","Rio de Janeiro/RJ ","22770-232")
  ,@("Rua Professor Henrique Costa - até 726 - lado par ","Pechincha ","Rio de Janeiro/RJ ","22770-233")
  ,@("Rua Professor Henrique Costa - de 687 ao fim - lado ímpar ","Pechincha ","Rio de Janeiro/RJ ","22770-234")
  ,@("Rua Professor Henrique Costa - de 728 ao fim - lado par ","Pechincha ","Rio de Janeiro/RJ ","22770-235")
  ,@("Rua Professor Mário de Vasconcelos ","Pechincha ","Rio de Janeiro/RJ ","22770-090")
  ,@("Rua Professor Rocha Lagoa ","Cidade de Deus ","Rio de Janeiro/RJ ","22770-390")
  ,@("Rua Professor Waldemar Berardinelli ","Cidade de Deus ","Rio de Janeiro/RJ ","22770-400")
  ,@("Rua Rebouças ","Cidade de Deus ","Rio de Janeiro/RJ ","22770-410")
  ,@("Rua Retiro dos Artistas - até 919 - lado ímpar ","Pechincha ","Rio de Janeiro/RJ ","22770-102")
  ,@("Rua Retiro dos Artistas - até 930 - lado par ","Pechincha ","Rio de Janeiro/RJ ","22770-103")
  ,@("Rua Retiro dos Artistas - de 921 ao fim - lado ímpar ","Pechincha ","Rio de Janeiro/RJ ","22770-104")
  ,@("Rua Retiro dos Artistas - de 932 ao fim - lado par ","Pechincha ","Rio de Janeiro/RJ ","22770-105")
  ,@("Rua Samuel das Neves ","Pechincha ","Rio de Janeiro/RJ ","22770-110")
  ,@("Rua São Boneto ","Pechincha ","Rio de Janeiro/RJ ","22770-430")
  ,@("Rua São Deodato ","Pechincha ","Rio de Janeiro/RJ ","22770-240")
  ,@("Rua Solar ","Cidade de Deus ","Rio de Janeiro/RJ ","22770-651")
  ,@("Rua Soldado Dirceu de Almeida ","Pechincha ","Rio de Janeiro/RJ ","22770-120")
  ,@("Rua Soldado Francisco de Souza ","Pechincha ","Rio de Janeiro/RJ ","22770-155")
  ,@("Rua Soldado Genésio Correia ","Taquara ","Rio de Janeiro/RJ ","22770-440")
  ,@("Rua Soldado Hilário Zanesco ","Taquara ","Rio de Janeiro/RJ ","22770-450")
  ,@("Rua Soldado João da Silva ","Taquara ","Rio de Janeiro/RJ ","22770-460")
  ,@("Rua Soldado José Solano ","Taquara ","Rio de Janeiro/RJ ","22770-470")
  ,@("Rua Tenente José Jerônimo de Mesquita ","Taquara ","Rio de Janeiro/RJ ","22770-250")
  ,@("Rua Vila Formosa ","Pechincha ","Rio de Janeiro/RJ ","22770-084")
  ,@("Rua Vila Lumiar ","Pechincha ","Rio de Janeiro/RJ ","22770-082")
  ,@("Rua Waldemar Loureiro ","Pechincha ","Rio de Janeiro/RJ ","22770-480")
  ,@("Rua Waldemar Rodrigues Martins ","Pechincha ","Rio de Janeiro/RJ ","22770-106")
)

$startRow = 51
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.StandardWidth = 20.115
